$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F18").Value = 150
$ws.Range("G18").Value = 17290.5
$ws.Range("F26").Value = 18
$ws.Range("G26").Value = 460.98
$ws.Range("F27").Value = 39
$ws.Range("G27").Value = 1398.54
$ws.Range("F28").Value = 44
$ws.Range("G28").Value = 1352.56
$ws.Range("F30").Value = 126
$ws.Range("G30").Value = 3550.68
$ws.Range("F31").Value = 34
$ws.Range("G31").Value = 905.76
$ws.Range("F32").Value = 29
$ws.Range("G32").Value = 1402.44
$ws.Range("B34").Value = 56203.74
$ws.Range("F36").Value = 79
$ws.Range("G36").Value = 15544.83
$ws.Range("F41").Value = 203
$ws.Range("G41").Value = 39156.67
$ws.Range("F52").Value = 24
$ws.Range("G52").Value = 1416
$ws.Range("F55").Value = 116
$ws.Range("G55").Value = 6468.16
$ws.Range("F58").Value = 61
$ws.Range("G58").Value = 4753.73
$ws.Range("F61").Value = 210
$ws.Range("G61").Value = 54753.3
$ws.Range("B66").Value = 191751.09
$ws.Range("F114").Value = 20
$ws.Range("G114").Value = 934.8
$ws.Range("F115").Value = 9
$ws.Range("G115").Value = 512.1
$ws.Range("B123").Value = 70421.47
$ws.Range("B126").Value = 65258
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("B127").Value = 64196
$ws.Range("F127").Value = 1
$ws.Range("G127").Value = 32143.58
$ws.Range("F143").Value = 41
$ws.Range("G143").Value = 7580.49
$ws.Range("B147").Value = 20540.89
$ws.Range("F151").Value = 25
$ws.Range("G151").Value = 3334
$ws.Range("B155").Value = 35280.95
$ws.Range("F172").Value = 94
$ws.Range("G172").Value = 5970.88
$ws.Range("F173").Value = 40
$ws.Range("G173").Value = 3144.4
$ws.Range("F184").Value = 47
$ws.Range("G184").Value = 3854
$ws.Range("B193").Value = 62025.26
$ws.Range("F206").Value = 60
$ws.Range("G206").Value = 3888
$ws.Range("B208").Value = 3888
$ws.Range("F212").Value = 58
$ws.Range("G212").Value = 5167.22
$ws.Range("F213").Value = 201
$ws.Range("G213").Value = 25462.68
$ws.Range("F216").Value = 70
$ws.Range("G216").Value = 5201
$ws.Range("B218").Value = 72203.67
$ws.Range("F222").Value = 637
$ws.Range("G222").Value = 11784.5
$ws.Range("B229").Value = 20965.9
$ws.Range("F252").Value = 69
$ws.Range("G252").Value = 6151.35
$ws.Range("F255").Value = 9
$ws.Range("G255").Value = 2842.2
$ws.Range("F263").Value = 7
$ws.Range("G263").Value = 725.9
$ws.Range("F268").Value = 9
$ws.Range("G268").Value = 1144.89
$ws.Range("F269").Value = 9
$ws.Range("G269").Value = 771.48
$ws.Range("F278").Value = 29
$ws.Range("G278").Value = 3930.66
$ws.Range("B290").Value = 66194
$ws.Range("C290").Value = 'HIM-Total Care Baby Pants Diapers-M-9s'
$ws.Range("F290").Value = 22
$ws.Range("G290").Value = 1884.96
$ws.Range("B291").Value = 64983
$ws.Range("C291").Value = 'HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S'
$ws.Range("F291").Value = 6
$ws.Range("G291").Value = 514.08
$ws.Range("B292").Value = 64985
$ws.Range("C292").Value = 'HIM-TOTAL CARE BABY PANTS DRAPERS-XL-9S'
$ws.Range("F292").Value = 12
$ws.Range("G292").Value = 1052.4
$ws.Range("B293").Value = 66196
$ws.Range("C293").Value = 'HIM-Total Care Baby Pants Drapers-Xl-9S'
$ws.Range("F293").Value = 6
$ws.Range("G293").Value = 526.2
$ws.Range("B295").Value = 108575.1
$ws.Range("B297").Value = 63565
$ws.Range("E297").Value = 109.19
$ws.Range("F297").Value = 60
$ws.Range("G297").Value = 6162.6
$ws.Range("B298").Value = 61610
$ws.Range("E298").Value = 122.71
$ws.Range("F298").Value = -58
$ws.Range("G298").Value = -5957.18
$ws.Range("F305").Value = 36
$ws.Range("G305").Value = 5194.08
$ws.Range("B306").Value = 57802
$ws.Range("E306").Value = 162.71
$ws.Range("F306").Value = -79
$ws.Range("G306").Value = -11334.92
$ws.Range("B307").Value = 63531
$ws.Range("E307").Value = 152.53
$ws.Range("F307").Value = 26
$ws.Range("G307").Value = 3730.48
$ws.Range("B328").Value = -17088.6
$ws.Range("F361").Value = 214
$ws.Range("G361").Value = 30086.26
$ws.Range("B363").Value = 70296.00999999999
$ws.Range("F366").Value = 52
$ws.Range("G366").Value = 2877.16
$ws.Range("F367").Value = 124
$ws.Range("G367").Value = 7526.8
$ws.Range("F370").Value = 196
$ws.Range("G370").Value = 32534.04
$ws.Range("F371").Value = 64
$ws.Range("G371").Value = 9617.280000000001
$ws.Range("B372").Value = 55535.6
$ws.Range("B381").Value = 58047
$ws.Range("D381").Value = 105.54
$ws.Range("E381").Value = 126.1
$ws.Range("F381").Value = 32
$ws.Range("G381").Value = 3377.28
$ws.Range("B382").Value = 47097
$ws.Range("D382").Value = 112.28
$ws.Range("E382").Value = 134.16
$ws.Range("F382").Value = 15
$ws.Range("G382").Value = 1684.2
$ws.Range("F402").Value = 37
$ws.Range("G402").Value = 1269.47
$ws.Range("B417").Value = 164472.25
$ws.Range("F432").Value = 77
$ws.Range("G432").Value = 3727.57
$ws.Range("B438").Value = 23536.34
$ws.Range("F454").Value = 63
$ws.Range("G454").Value = 17826.48
$ws.Range("B458").Value = 89389.38
$ws.Range("F472").Value = 9
$ws.Range("G472").Value = 11410.38
$ws.Range("B476").Value = 43932.69
$ws.Range("F511").Value = 207
$ws.Range("G511").Value = 20673.09
$ws.Range("B525").Value = 115948.56
$ws.Range("F528").Value = 251
$ws.Range("G528").Value = 3980.86
$ws.Range("B535").Value = 22027.07
$ws.Range("F551").Value = 9
$ws.Range("G551").Value = 6706.71
$ws.Range("B556").Value = 41164.62
$ws.Range("F560").Value = 22
$ws.Range("G560").Value = 1768.36
$ws.Range("B561").Value = 24335.87
$ws.Range("F563").Value = 15
$ws.Range("G563").Value = 2802.9
$ws.Range("B573").Value = 15963.73
$ws.Range("F605").Value = 170
$ws.Range("G605").Value = 22627
$ws.Range("B607").Value = 23032.03
$ws.Range("F609").Value = 7
$ws.Range("G609").Value = 761.67
$ws.Range("F614").Value = 77
$ws.Range("G614").Value = 11171.16
$ws.Range("F623").Value = 77
$ws.Range("G623").Value = 39619.58
$ws.Range("F625").Value = 312
$ws.Range("G625").Value = 11490.96
$ws.Range("F626").Value = 9
$ws.Range("G626").Value = 424.89
$ws.Range("B628").Value = 201717.76
$ws.Range("F660").Value = 47
$ws.Range("G660").Value = 1397.78
$ws.Range("F662").Value = 35
$ws.Range("G662").Value = 2810.85
$ws.Range("B668").Value = 11105.06
$ws.Range("F703").Value = 5
$ws.Range("G703").Value = 2426.35
$ws.Range("F706").Value = 111
$ws.Range("G706").Value = 4350.09
$ws.Range("F711").Value = 11
$ws.Range("G711").Value = 5872.35
$ws.Range("B713").Value = 61863.29
$ws.Range("B718").Value = 2477456.24
$ws.Range("B719").Value = 2477456.24
